{"js": "// Update the date line and the 25 \"a\u00f7b=\" division prompts in the practice\n// table. Each table cell is addressed by its row/column position (rather\n// than a blind global text search) because a couple of the original\n// prompts (e.g. \"67\u00f78=\") repeat verbatim but must become different values\n// depending on where they sit in the table.\n\n// 1) Date heading paragraph: \"2025-06-25 Wednesday\" -> \"2025-06-26 Thursday\"\nconst titleResults = context.document.body.search(\"2025-06-25 Wednesday\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"2025-06-26 Thursday\", Word.InsertLocation.replace);\n}\n\n// 2) The division prompts live in the first row of every 4-row block in the\n// single table on the page. Address them as (rowIndex, columnIndex) pairs\n// so duplicate prompt text (\"67\u00f78=\" appears twice) resolves unambiguously.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Old value kept only as a sanity check / documentation; the new value is\n// what actually gets written, addressed purely by position.\nconst grid = [\n  { row: 0, col: 0, oldText: \"11\u00f78=\", newText: \"44\u00f77=\" },\n  { row: 0, col: 1, oldText: \"76\u00f78=\", newText: \"92\u00f75=\" },\n  { row: 0, col: 2, oldText: \"51\u00f76=\", newText: \"27\u00f73=\" },\n  { row: 0, col: 3, oldText: \"89\u00f79=\", newText: \"83\u00f76=\" },\n  { row: 0, col: 4, oldText: \"17\u00f76=\", newText: \"69\u00f78=\" },\n\n  { row: 4, col: 0, oldText: \"30\u00f78=\", newText: \"63\u00f72=\" },\n  { row: 4, col: 1, oldText: \"39\u00f79=\", newText: \"14\u00f75=\" },\n  { row: 4, col: 2, oldText: \"67\u00f78=\", newText: \"86\u00f76=\" },\n  { row: 4, col: 3, oldText: \"60\u00f78=\", newText: \"46\u00f76=\" },\n  { row: 4, col: 4, oldText: \"48\u00f73=\", newText: \"51\u00f79=\" },\n\n  { row: 8, col: 0, oldText: \"78\u00f77=\", newText: \"23\u00f72=\" },\n  { row: 8, col: 1, oldText: \"61\u00f76=\", newText: \"58\u00f73=\" },\n  { row: 8, col: 2, oldText: \"55\u00f77=\", newText: \"20\u00f77=\" },\n  { row: 8, col: 3, oldText: \"52\u00f72=\", newText: \"87\u00f79=\" },\n  { row: 8, col: 4, oldText: \"81\u00f76=\", newText: \"47\u00f73=\" },\n\n  { row: 12, col: 0, oldText: \"58\u00f75=\", newText: \"97\u00f77=\" },\n  { row: 12, col: 1, oldText: \"19\u00f78=\", newText: \"37\u00f78=\" },\n  { row: 12, col: 2, oldText: \"92\u00f79=\", newText: \"41\u00f78=\" },\n  { row: 12, col: 3, oldText: \"98\u00f79=\", newText: \"74\u00f73=\" },\n  { row: 12, col: 4, oldText: \"59\u00f79=\", newText: \"10\u00f73=\" },\n\n  { row: 16, col: 0, oldText: \"75\u00f77=\", newText: \"91\u00f72=\" },\n  { row: 16, col: 1, oldText: \"55\u00f75=\", newText: \"56\u00f72=\" },\n  { row: 16, col: 2, oldText: \"51\u00f75=\", newText: \"60\u00f75=\" },\n  { row: 16, col: 3, oldText: \"33\u00f72=\", newText: \"61\u00f75=\" },\n  { row: 16, col: 4, oldText: \"67\u00f78=\", newText: \"57\u00f78=\" },\n];\n\n// Load all needed cell collections up front.\nconst neededRowIdx = [...new Set(grid.map((g) => g.row))];\nconst rowCells = {};\nfor (const ri of neededRowIdx) {\n  const cells = rows.items[ri].cells;\n  cells.load(\"items\");\n  rowCells[ri] = cells;\n}\nawait context.sync();\n\n// For each target cell, do an in-cell search-and-replace so the original\n// run formatting (font, size, paragraph alignment) is preserved.\nfor (const entry of grid) {\n  const cell = rowCells[entry.row].items[entry.col];\n  const cellResults = cell.body.search(entry.oldText, { matchCase: true });\n  cellResults.load(\"items\");\n  await context.sync();\n  if (cellResults.items.length > 0) {\n    cellResults.items[0].insertText(entry.newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: if for some reason the expected old text isn't present,\n    // replace the whole cell body text directly.\n    cell.body.insertText(entry.newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"a\u00f7b=\" division prompts in the practice\n# table. Each table cell is addressed by its (row, column) position (rather\n# than a blind document-wide Find/Replace) because a couple of the original\n# prompts (e.g. \"67\u00f78=\") repeat verbatim but must become different values\n# depending on where they sit in the table.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading: \"2025-06-25 Wednesday\" -> \"2025-06-26 Thursday\"\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleRange.Find.Execute(\"2025-06-25 Wednesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-06-26 Thursday\", 2)\n\n# 2) Division prompts, addressed by (row, column) in the single table.\n#    Word table cell indices are 1-based; the prompts occupy the first row\n#    of every 4-row block (rows 1, 5, 9, 13, 17).\n$t = $d.Tables.Item(1)\n\n$grid = @(\n    @{ Row = 1;  Col = 1; New = \"44\u00f77=\" },\n    @{ Row = 1;  Col = 2; New = \"92\u00f75=\" },\n    @{ Row = 1;  Col = 3; New = \"27\u00f73=\" },\n    @{ Row = 1;  Col = 4; New = \"83\u00f76=\" },\n    @{ Row = 1;  Col = 5; New = \"69\u00f78=\" },\n\n    @{ Row = 5;  Col = 1; New = \"63\u00f72=\" },\n    @{ Row = 5;  Col = 2; New = \"14\u00f75=\" },\n    @{ Row = 5;  Col = 3; New = \"86\u00f76=\" },\n    @{ Row = 5;  Col = 4; New = \"46\u00f76=\" },\n    @{ Row = 5;  Col = 5; New = \"51\u00f79=\" },\n\n    @{ Row = 9;  Col = 1; New = \"23\u00f72=\" },\n    @{ Row = 9;  Col = 2; New = \"58\u00f73=\" },\n    @{ Row = 9;  Col = 3; New = \"20\u00f77=\" },\n    @{ Row = 9;  Col = 4; New = \"87\u00f79=\" },\n    @{ Row = 9;  Col = 5; New = \"47\u00f73=\" },\n\n    @{ Row = 13; Col = 1; New = \"97\u00f77=\" },\n    @{ Row = 13; Col = 2; New = \"37\u00f78=\" },\n    @{ Row = 13; Col = 3; New = \"41\u00f78=\" },\n    @{ Row = 13; Col = 4; New = \"74\u00f73=\" },\n    @{ Row = 13; Col = 5; New = \"10\u00f73=\" },\n\n    @{ Row = 17; Col = 1; New = \"91\u00f72=\" },\n    @{ Row = 17; Col = 2; New = \"56\u00f72=\" },\n    @{ Row = 17; Col = 3; New = \"60\u00f75=\" },\n    @{ Row = 17; Col = 4; New = \"61\u00f75=\" },\n    @{ Row = 17; Col = 5; New = \"57\u00f78=\" }\n)\n\nforeach ($entry in $grid) {\n    $cell = $t.Cell($entry.Row, $entry.Col)\n    $r = $cell.Range\n    # Trim the trailing cell-end mark so we only overwrite the visible text,\n    # leaving the run's formatting (font/size) and paragraph alignment intact.\n    $r.End = $r.End - 1\n    $r.Text = $entry.New\n}\n"}
